$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("site_metrics")
$ws1.Range("O15").Value = 0.005779341191785079

$ws1.Range("O16").Value = 0.01789475208351017

$ws1.Range("O17").Value = 0.005618578486865572

$ws1.Range("AK17").Value = $true

$ws1.Range("AK20").Value = $true

$ws1.Range("AK21").Value = $true

$ws1.Range("O22").Value = 0.04698714589605456

$ws1.Range("O23").Value = 0.04900963588476684

$ws1.Range("AK23").Value = $true

$ws1.Range("O28").Value = 0.002786193505485551

$ws1.Range("AK42").Value = $true

$ws1.Range("O43").Value = 0.04307847649894682

$ws1.Range("O57").Value = 0.2156407608830078

$ws1.Range("AK58").Value = $true

$ws1.Range("AK62").Value = $true

$ws1.Range("O68").Value = 0.04734500304124714

$ws1.Range("O73").Value = 0.0212974678578754

$ws1.Range("O78").Value = 0.00330623289296553

$ws1.Range("AK78").Value = $true

$ws1.Range("O80").Value = 0.002732410185299382

$ws1.Range("O81").Value = 0.004649956386151363

$ws1.Range("O82").Value = 0.04609438884185086

$ws1.Range("O83").Value = 0.0481278083251732

$ws1.Range("O86").Value = 0.0199608926356963

$ws1.Range("AK88").Value = $true

$ws1.Range("AK91").Value = $true

$ws1.Range("O93").Value = 0.07324438939577617

$ws1.Range("O94").Value = 0.007163830006606227

$ws1.Range("AK96").Value = $true

$ws1.Range("AK99").Value = $true

$ws1.Range("O101").Value = 0.0481647123529429

$ws1.Range("O103").Value = 0.04403517620970954

$ws1.Range("O104").Value = 0.05533164632087837

$ws1.Range("AK119").Value = $true

$ws1.Range("O121").Value = 0.01162218826395471

$ws1.Range("AK126").Value = $true

$ws1.Range("AK128").Value = $true

$ws1.Range("O131").Value = 0.1322786678075588

$ws1.Range("AK132").Value = $true

$ws1.Range("O133").Value = 0.02855739494976635

$ws1.Range("AK133").Value = $true

$ws1.Range("AK135").Value = $true

$ws1.Range("AK136").Value = $true

$ws1.Range("O142").Value = 0.006994705988248208

$ws3 = $wb.Worksheets.Item("mk_duration")
$ws3.Range("M4").Value = 0.09124320211532533
$ws3.Range("N4").Value = 1.688875965185925
$ws3.Range("O4").Value = 0.2192118226600985
$ws3.Range("P4").Value = 89
$ws3.Range("Q4").Value = 2715
$ws3.Range("R4").Value = 0.2739583333333333
$ws3.Range("S4").Value = 2.164583333333334

$ws3.Range("M26").Value = 0.9046033291427005
$ws3.Range("N26").Value = 0.1198482881916595
$ws3.Range("O26").Value = 0.01724137931034483
$ws3.Range("P26").Value = 7
$ws3.Range("Q26").Value = 2506.333333333333
$ws3.Range("S26").Value = 4

$ws3.Range("M29").Value = 0.05990572743327505
$ws3.Range("N29").Value = 1.881486813157878
$ws3.Range("O29").Value = 0.2438423645320197
$ws3.Range("P29").Value = 99
$ws3.Range("Q29").Value = 2713
$ws3.Range("R29").Value = 0.1696969696969697
$ws3.Range("S29").Value = -0.375757575757576

$ws3.Range("M40").Value = 0.9334856116416619
$ws3.Range("N40").Value = 0.0834602139578355
$ws3.Range("O40").Value = 0.01424501424501425
$ws3.Range("P40").Value = 5
$ws3.Range("Q40").Value = 2297
$ws3.Range("R40").Value = 0.01515151515151518
$ws3.Range("S40").Value = 4.088744588744588

$ws3.Range("M43").Value = 0.2886296994062758
$ws3.Range("N43").Value = 1.061132478955467
$ws3.Range("O43").Value = 0.1428571428571428
$ws3.Range("P43").Value = 54
$ws3.Range("Q43").Value = 2494.666666666667
$ws3.Range("R43").Value = 0.096875
$ws3.Range("S43").Value = 2.525520833333333

$ws3.Range("K46").Value = "no trend"
$ws3.Range("L46").Value = $false
$ws3.Range("M46").Value = 0.577274383745257
$ws3.Range("N46").Value = -0.5573704017131537
$ws3.Range("O46").Value = -0.1029411764705882
$ws3.Range("P46").Value = -14
$ws3.Range("Q46").Value = 544
$ws3.Range("S46").Value = 8

$ws3.Range("K72").Value = "no trend"
$ws3.Range("L72").Value = $false
$ws3.Range("M72").Value = 0.08384031020703486
$ws3.Range("N72").Value = 1.728825615270013
$ws3.Range("O72").Value = 0.2597402597402597
$ws3.Range("P72").Value = 60
$ws3.Range("Q72").Value = 1164.666666666667
$ws3.Range("R72").Value = 0.1666666666666667
$ws3.Range("S72").Value = 0.08333333333333348

$ws3.Range("K122").Value = "no trend"
$ws3.Range("L122").Value = $false
$ws3.Range("M122").Value = 0.1350836263468835
$ws3.Range("N122").Value = 1.49435205382276
$ws3.Range("O122").Value = 0.225296442687747
$ws3.Range("P122").Value = 57
$ws3.Range("Q122").Value = 1404.333333333333
$ws3.Range("R122").Value = 0.3015873015873016
$ws3.Range("S122").Value = 3.015873015873016

$ws3.Range("K127").Value = "no trend"
$ws3.Range("L127").Value = $false
$ws3.Range("M127").Value = 0.05112272749065738
$ws3.Range("N127").Value = 1.950448151275792
$ws3.Range("O127").Value = 0.3529411764705883
$ws3.Range("P127").Value = 48
$ws3.Range("Q127").Value = 580.6666666666666
$ws3.Range("R127").Value = 0.7895833333333333
$ws3.Range("S127").Value = 1.350000000000001

$ws4 = $wb.Worksheets.Item("mk_intra_annual")
$ws4.Range("M4").Value = 0.1692354188501368
$ws4.Range("N4").Value = 1.374664703360094
$ws4.Range("O4").Value = 0.1748768472906404
$ws4.Range("P4").Value = 71
$ws4.Range("Q4").Value = 2593

$ws4.Range("M26").Value = 0.713284142388857
$ws4.Range("N26").Value = 0.3674489795938048
$ws4.Range("O26").Value = 0.04679802955665024
$ws4.Range("P26").Value = 19
$ws4.Range("Q26").Value = 2399.666666666667
$ws4.Range("S26").Value = 1

$ws4.Range("M29").Value = 0.4170480732510811
$ws4.Range("N29").Value = 0.8115531676178501
$ws4.Range("O29").Value = 0.1059113300492611
$ws4.Range("P29").Value = 43
$ws4.Range("Q29").Value = 2678.333333333333
$ws4.Range("S29").Value = 2

$ws4.Range("M40").Value = 1
$ws4.Range("N40").Value = 0
$ws4.Range("O40").Value = 0
$ws4.Range("P40").Value = 0
$ws4.Range("Q40").Value = 2090

$ws4.Range("K43").Value = "no trend"
$ws4.Range("L43").Value = $false
$ws4.Range("M43").Value = 0.05668784899732704
$ws4.Range("N43").Value = 1.905709939635053
$ws4.Range("O43").Value = 0.2513227513227513
$ws4.Range("P43").Value = 95
$ws4.Range("Q43").Value = 2433
$ws4.Range("R43").Value = 0.08012820512820512
$ws4.Range("S43").Value = -0.08173076923076916

$ws4.Range("K46").Value = "no trend"
$ws4.Range("L46").Value = $false
$ws4.Range("M46").Value = 0.6151205041013237
$ws4.Range("N46").Value = -0.5027777991522042
$ws4.Range("O46").Value = -0.08823529411764706
$ws4.Range("P46").Value = -12
$ws4.Range("Q46").Value = 478.6666666666667
$ws4.Range("S46").Value = 1

$ws4.Range("M72").Value = 0.7902291355909628
$ws4.Range("N72").Value = 0.2660130798453453
$ws4.Range("O72").Value = 0.04329004329004329
$ws4.Range("P72").Value = 10
$ws4.Range("Q72").Value = 1144.666666666667
$ws4.Range("S72").Value = 1

$ws4.Range("M122").Value = 0.0136468100021383
$ws4.Range("N122").Value = 2.466428160553326
$ws4.Range("O122").Value = 0.3517786561264822
$ws4.Range("P122").Value = 89
$ws4.Range("Q122").Value = 1273
$ws4.Range("R122").Value = 0.06666666666666667
$ws4.Range("S122").Value = 0.2666666666666667

$ws4.Range("M127").Value = 0.9323851505432275
$ws4.Range("N127").Value = -0.08484432973359157
$ws4.Range("O127").Value = -0.02205882352941177
$ws4.Range("P127").Value = -3
$ws4.Range("Q127").Value = 555.6666666666666
$ws4.Range("S127").Value = 2
